$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B values (row labels use col B still, pre-delete)
$ws.Cells.Item(2, 2).Value = 6422.554661072401
$ws.Cells.Item(3, 2).Value = 4094.559754097033
$ws.Cells.Item(4, 2).Value = 5840.972492930345
$ws.Cells.Item(5, 2).Value = 4503.490017348013
$ws.Cells.Item(6, 2).Value = 5527.399507626282
$ws.Cells.Item(7, 2).Value = 5322.457917588043
$ws.Cells.Item(8, 2).Value = 5497.207733982736
$ws.Cells.Item(9, 2).Value = 9320.372100309036
$ws.Cells.Item(10, 2).Value = 6723.246708300027
$ws.Cells.Item(11, 2).Value = 9302.877988982826
$ws.Cells.Item(12, 2).Value = 11112.98756336252
$ws.Cells.Item(13, 2).Value = 4657.974863404017
$ws.Cells.Item(14, 2).Value = 4486.473016585671
$ws.Cells.Item(15, 2).Value = 4622.675904585194
$ws.Cells.Item(16, 2).Value = 4869.069683910613
$ws.Cells.Item(17, 2).Value = 4869.069683910613
$ws.Cells.Item(18, 2).Value = 4869.069683910613
$ws.Cells.Item(19, 2).Value = 4869.069683910613
$ws.Cells.Item(20, 2).Value = 4292.357644096977
$ws.Cells.Item(21, 2).Value = 4292.357644096977
$ws.Cells.Item(22, 2).Value = 4869.069683910613
$ws.Cells.Item(23, 2).Value = 4869.069683910613
$ws.Cells.Item(24, 2).Value = 4869.069683910613
$ws.Cells.Item(25, 2).Value = 4869.069683910613
$ws.Cells.Item(26, 2).Value = 4869.069683910613
$ws.Cells.Item(27, 2).Value = 4292.357644096977
$ws.Cells.Item(28, 2).Value = 4869.069683910613
$ws.Cells.Item(29, 2).Value = 4869.069683910613
$ws.Cells.Item(30, 2).Value = 4292.357644096977
$ws.Cells.Item(31, 2).Value = 9662.041362311038
$ws.Cells.Item(32, 2).Value = 5740.240021630538

# Remove the "max" column (C); D (prediction) and E (rejection-f) shift left
$ws.Columns.Item(3).Delete()

